$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency list refresh (GitHub Actions scheduled run): refresh the
# "Price" (D) and "Volume(1h)" (E) columns for every coin row (2-51).
#
# The Price column stores plain text (e.g. "1.001", "290.78") rather than
# numbers, so the whole column is switched to a Text number format before
# the writes (keeping numeric-looking quotes such as "1.001" intact instead
# of being parsed into floating-point numbers) and then restored to the
# workbook default style afterwards so formatting is left untouched.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '22.399.82'
$ws.Range("E2").Value = '  -4.67%  '
$ws.Range("D3").Value = '1.569.66'
$ws.Range("E3").Value = '  -4.74%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '290.78'
$ws.Range("E6").Value = '  -2.95%  '
$ws.Range("D7").Value = '0.3689'
$ws.Range("E7").Value = '  -2.78%  '
$ws.Range("D8").Value = '49.64'
$ws.Range("E8").Value = '  -1.00%  '
$ws.Range("D9").Value = '0.3364'
$ws.Range("E9").Value = '  -5.77%  '
$ws.Range("E10").Value = '  -4.49%  '
$ws.Range("D11").Value = '0.07561'
$ws.Range("E11").Value = '  -6.53%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").Value = '21.04'
$ws.Range("E13").Value = '  -4.42%  '
$ws.Range("D14").Value = '6.048'
$ws.Range("E14").Value = '  -5.49%  '
$ws.Range("D15").Value = '6.846'
$ws.Range("E15").Value = '  -7.22%  '
$ws.Range("D16").Value = '0.00001142'
$ws.Range("E16").Value = '  -4.50%  '
$ws.Range("D17").Value = '1.570.16'
$ws.Range("E17").Value = '  -4.74%  '
$ws.Range("D18").Value = '89.19'
$ws.Range("E18").Value = '  -8.11%  '
$ws.Range("D19").Value = '0.06666'
$ws.Range("E19").Value = '  -4.30%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = '6.227'
$ws.Range("E21").Value = '  -7.79%  '
$ws.Range("D22").Value = '16.33'
$ws.Range("E22").Value = '  -6.00%  '
$ws.Range("D23").Value = '11.96'
$ws.Range("E23").Value = '  -3.64%  '
$ws.Range("D24").Value = '22.405.58'
$ws.Range("E24").Value = '  -4.73%  '
$ws.Range("D25").Value = '2.396'
$ws.Range("E25").Value = '  -4.13%  '
$ws.Range("D26").Value = '2.952'
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("D27").Value = '19.83'
$ws.Range("E27").Value = '  -5.19%  '
$ws.Range("D28").Value = '146.28'
$ws.Range("E28").Value = '  -4.35%  '
$ws.Range("D29").Value = '4.931'
$ws.Range("E29").Value = '  -5.46%  '
$ws.Range("D30").Value = '125.05'
$ws.Range("E30").Value = '  -5.73%  '
$ws.Range("D31").Value = '1.745.51'
$ws.Range("E31").Value = '  -4.47%  '
$ws.Range("D32").Value = '6.263'
$ws.Range("E32").Value = '  -9.28%  '
$ws.Range("D33").Value = '1.970'
$ws.Range("E33").Value = '  -7.82%  '
$ws.Range("D34").Value = '0.9755'
$ws.Range("E34").Value = '  -4.23%  '
$ws.Range("D35").Value = '10.36'
$ws.Range("E35").Value = '  -12.45%  '
$ws.Range("D36").Value = '0.08427'
$ws.Range("E36").Value = '  -3.60%  '
$ws.Range("D37").Value = '0.02523'
$ws.Range("E37").Value = '  -7.58%  '
$ws.Range("D38").Value = '0.2300'
$ws.Range("E38").Value = '  -5.45%  '
$ws.Range("D39").Value = '0.06504'
$ws.Range("E39").Value = '  -4.27%  '
$ws.Range("D40").Value = '5.480'
$ws.Range("E40").Value = '  -7.68%  '
$ws.Range("D41").Value = '11.77'
$ws.Range("E41").Value = '  -10.35%  '
$ws.Range("D42").Value = '1.244'
$ws.Range("E42").Value = '  -5.62%  '
$ws.Range("D43").Value = '0.6381'
$ws.Range("E43").Value = '  -7.45%  '
$ws.Range("D44").Value = '14.47'
$ws.Range("E44").Value = '  -6.22%  '
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("D46").Value = '0.5989'
$ws.Range("E46").Value = '  -6.51%  '
$ws.Range("D47").Value = '3.767'
$ws.Range("E47").Value = '  -4.04%  '
$ws.Range("E48").Value = '  -6.91%  '
$ws.Range("D49").Value = '121.42'
$ws.Range("E49").Value = '  -4.88%  '
$ws.Range("D50").Value = '0.07264'
$ws.Range("E50").Value = '  -6.21%  '
$ws.Range("D51").Value = '1.185'
$ws.Range("E51").Value = '  +0.23%  '

$priceRange.Style = "Normal"
